# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Rows map to events whose interest counts increased:
#   row 5  : 492  -> 493
#   row 6  : 1431 -> 1440
#   row 7  : 723  -> 747
#   row 9  : 190  -> 191
#   row 10 : 137  -> 138
#   row 11 : 188  -> 189
#   row 13 : 162  -> 163

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 493
    6  = 1440
    7  = 747
    9  = 191
    10 = 138
    11 = 189
    13 = 163
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
